# Add a new worksheet "Yearly demand" at the end of the workbook with
# the hourly demand data for three representative days.

$wb = $excel.ActiveWorkbook

# Reference sheet that already has the exact same layout/style (header row
# B1:Y1, row-label column A2:A4, page margins, outline/page-setup props)
# so we duplicate it rather than building the sheet from scratch - this
# keeps all formatting/style references identical to the rest of the
# workbook instead of minting new ones.
$template = $wb.Worksheets.Item("DG Dispatch")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Yearly demand"

# Header row (hour numbers 0-23) in B1:Y1
for ($i = 0; $i -le 23; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $i
}

# Row label column (A2:A4) values 0,1,2
for ($r = 0; $r -le 2; $r++) {
    $ws.Cells.Item(2 + $r, 1).Value = $r
}

$data = @(
    @(-32.5, -19.5, -13, -13, -13, 142.5, 291.5, 327, 388.5, 502, 596, 670.5, 745, 651, 576.5, 502, 320.5, 139, 32, -117, -97.5, -78, -52, -39),
    @(-32.5, -19.5, -13, 0, 0, -19.5, 0, 324, 486, 648, 729, 751.5, 583, 567, 333.5, 340, 243, 57.99999999999999, -130, 0, 0, -78, 0, -39),
    @(-32.5, -19.5, 0, 0, 0, -19.5, 0, 0, 81, 324, 567, 589.5, 648, 567, 324, 162, 81, 0, -130, 0, 0, 0, 0, -39)
)

for ($r = 0; $r -le 2; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -le 23; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value = $rowData[$c]
    }
}

# Restore the originally active sheet/selection (adding/copying a sheet
# shifts focus to it in Excel) so the workbook-level view state is left
# untouched, matching the source edit.
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()
